# Applies the daily cryptos.xlsx data refresh (commit: "Updated cryptos list ...
# with GitHub Actions"): refreshed Price + Volume(1h) figures for every coin row,
# plus a coin-ranking reshuffle affecting rows 36-39 (Coin name + Link + Price +
# Volume(1h) all move to reflect the new ordering).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numbers as literal text, e.g. "63.078.20" or
# "6.40", using '.' as a thousands separator and sometimes trailing zeros that
# must be preserved exactly. Plain `.Value = "6.40"` lets Excel "smart type"
# the text into a real number (dropping the trailing zero, turning
# "63.078.20" shaped values into dates/numbers, etc). Prefixing with an
# apostrophe forces text entry like a user typing it in the UI; ClearFormats()
# then drops the resulting quote-prefix cell style so the cell ends up with no
# style at all, matching every other (untouched) cell in this column.
function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "63.123.93"
$ws.Range("E2").Value = "  +1.76%  "
Set-TextValue $ws.Range("D3") "3.459.54"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue $ws.Range("D5") "580.51"
$ws.Range("E5").Value = "  +0.37%  "
Set-TextValue $ws.Range("D6") "147.59"
$ws.Range("E6").Value = "  +1.80%  "
Set-TextValue $ws.Range("D7") "3.460.18"
$ws.Range("E7").Value = "  +1.24%  "
$ws.Range("E8").Value = "  -0.07%  "
Set-TextValue $ws.Range("D9") "0.479"
$ws.Range("E9").Value = "  +0.79%  "
Set-TextValue $ws.Range("D10") "7.83"
$ws.Range("E10").Value = "  +3.00%  "
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("E12").Value = "  +4.93%  "
Set-TextValue $ws.Range("D13") "4.053.19"
$ws.Range("E13").Value = "  +1.26%  "
Set-TextValue $ws.Range("D14") "29.24"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("E15").Value = "  +2.41%  "
Set-TextValue $ws.Range("D16") "3.461.71"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("E17").Value = "  +0.89%  "
Set-TextValue $ws.Range("D18") "63.123.45"
$ws.Range("E18").Value = "  +1.82%  "
Set-TextValue $ws.Range("D19") "6.42"
$ws.Range("E19").Value = "  +3.94%  "
Set-TextValue $ws.Range("D20") "14.45"
$ws.Range("E20").Value = "  +2.59%  "
Set-TextValue $ws.Range("D21") "9.27"
$ws.Range("E21").Value = "  +0.95%  "
Set-TextValue $ws.Range("D22") "386.66"
$ws.Range("E22").Value = "  -1.11%  "
Set-TextValue $ws.Range("D23") "0.563"
$ws.Range("E23").Value = "  +1.33%  "
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("E25").Value = "  -0.01%  "
Set-TextValue $ws.Range("D26") "3.602.89"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  -2.46%  "
Set-TextValue $ws.Range("D29") "7.64"
$ws.Range("E29").Value = "  +2.13%  "
Set-TextValue $ws.Range("D30") "1.00"
$ws.Range("E30").Value = "  +0.03%  "
Set-TextValue $ws.Range("D31") "8.16"
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("E34").Value = "  -4.08%  "
Set-TextValue $ws.Range("D35") "23.36"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D36") "5.32"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D37") "7.14"
$ws.Range("E37").Value = "  +2.24%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D38") "32.05"
$ws.Range("E38").Value = "  +11.33%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D39") "1.60"
$ws.Range("E39").Value = "  +3.31%  "
Set-TextValue $ws.Range("D40") "168.26"
$ws.Range("E40").Value = "  +0.14%  "
Set-TextValue $ws.Range("D41") "3.496.75"
$ws.Range("E41").Value = "  +1.37%  "
Set-TextValue $ws.Range("D42") "0.0770"
$ws.Range("E42").Value = "  +1.96%  "
Set-TextValue $ws.Range("D43") "0.793"
$ws.Range("E43").Value = "  +0.80%  "
Set-TextValue $ws.Range("D44") "42.43"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("E45").Value = "  +3.31%  "
$ws.Range("E46").Value = "  +3.35%  "
$ws.Range("E47").Value = "  -1.69%  "
Set-TextValue $ws.Range("D48") "2.590.87"
$ws.Range("E48").Value = "  +3.48%  "
$ws.Range("E49").Value = "  +10.13%  "
Set-TextValue $ws.Range("D50") "6.80"
$ws.Range("E50").Value = "  +2.28%  "
Set-TextValue $ws.Range("D51") "22.94"
$ws.Range("E51").Value = "  +0.35%  "
